$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '24.392.67'
$ws.Range('D2').NumberFormat = "General"
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  -1.59%  '
$ws.Range('E2').NumberFormat = "General"
$ws.Range('E2').Style = "Normal"
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.653.92'
$ws.Range('D3').NumberFormat = "General"
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  -2.58%  '
$ws.Range('E3').NumberFormat = "General"
$ws.Range('E3').Style = "Normal"
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.002'
$ws.Range('D4').NumberFormat = "General"
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('E4').NumberFormat = "General"
$ws.Range('E4').Style = "Normal"
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '311.38'
$ws.Range('D5').NumberFormat = "General"
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  -1.20%  '
$ws.Range('E5').NumberFormat = "General"
$ws.Range('E5').Style = "Normal"
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  -0.06%  '
$ws.Range('E6').NumberFormat = "General"
$ws.Range('E6').Style = "Normal"
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.3914'
$ws.Range('D7').NumberFormat = "General"
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  -1.83%  '
$ws.Range('E7').NumberFormat = "General"
$ws.Range('E7').Style = "Normal"
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3907'
$ws.Range('D8').NumberFormat = "General"
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  -3.44%  '
$ws.Range('E8').NumberFormat = "General"
$ws.Range('E8').Style = "Normal"
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  +0.02%  '
$ws.Range('E9').NumberFormat = "General"
$ws.Range('E9').Style = "Normal"
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  -5.92%  '
$ws.Range('E10').NumberFormat = "General"
$ws.Range('E10').Style = "Normal"
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '50.10'
$ws.Range('D11').NumberFormat = "General"
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  -6.45%  '
$ws.Range('E11').NumberFormat = "General"
$ws.Range('E11').Style = "Normal"
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.08552'
$ws.Range('D12').NumberFormat = "General"
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  -2.86%  '
$ws.Range('E12').NumberFormat = "General"
$ws.Range('E12').Style = "Normal"
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '24.88'
$ws.Range('D13').NumberFormat = "General"
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  -5.41%  '
$ws.Range('E13').NumberFormat = "General"
$ws.Range('E13').Style = "Normal"
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '7.209'
$ws.Range('D14').NumberFormat = "General"
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  -4.51%  '
$ws.Range('E14').NumberFormat = "General"
$ws.Range('E14').Style = "Normal"
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.00001304'
$ws.Range('D15').NumberFormat = "General"
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  -2.95%  '
$ws.Range('E15').NumberFormat = "General"
$ws.Range('E15').Style = "Normal"
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '7.604'
$ws.Range('D16').NumberFormat = "General"
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  -4.73%  '
$ws.Range('E16').NumberFormat = "General"
$ws.Range('E16').Style = "Normal"
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '1.657.88'
$ws.Range('D17').NumberFormat = "General"
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  -3.96%  '
$ws.Range('E17').NumberFormat = "General"
$ws.Range('E17').Style = "Normal"
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '93.17'
$ws.Range('D18').NumberFormat = "General"
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  -2.45%  '
$ws.Range('E18').NumberFormat = "General"
$ws.Range('E18').Style = "Normal"
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  -3.15%  '
$ws.Range('E19').NumberFormat = "General"
$ws.Range('E19').Style = "Normal"
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '20.92'
$ws.Range('D20').NumberFormat = "General"
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  +0.10%  '
$ws.Range('E20').NumberFormat = "General"
$ws.Range('E20').Style = "Normal"
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '7.002'
$ws.Range('D21').NumberFormat = "General"
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  -4.66%  '
$ws.Range('E21').NumberFormat = "General"
$ws.Range('E21').Style = "Normal"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '1.001'
$ws.Range('D22').NumberFormat = "General"
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  -0.08%  '
$ws.Range('E22').NumberFormat = "General"
$ws.Range('E22').Style = "Normal"
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  -4.07%  '
$ws.Range('E23').NumberFormat = "General"
$ws.Range('E23').Style = "Normal"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '24.400.51'
$ws.Range('D24').NumberFormat = "General"
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  -1.50%  '
$ws.Range('E24').NumberFormat = "General"
$ws.Range('E24').Style = "Normal"
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.336'
$ws.Range('D25').NumberFormat = "General"
$ws.Range('D25').Style = "Normal"
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.784'
$ws.Range('D26').NumberFormat = "General"
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  -4.10%  '
$ws.Range('E26').NumberFormat = "General"
$ws.Range('E26').Style = "Normal"
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  -1.66%  '
$ws.Range('E27').NumberFormat = "General"
$ws.Range('E27').Style = "Normal"
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '158.67'
$ws.Range('D28').NumberFormat = "General"
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  -1.94%  '
$ws.Range('E28').NumberFormat = "General"
$ws.Range('E28').Style = "Normal"
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '5.727'
$ws.Range('D29').NumberFormat = "General"
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  -6.97%  '
$ws.Range('E29').NumberFormat = "General"
$ws.Range('E29').Style = "Normal"
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '145.37'
$ws.Range('D30').NumberFormat = "General"
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  +0.72%  '
$ws.Range('E30').NumberFormat = "General"
$ws.Range('E30').Style = "Normal"
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '8.212'
$ws.Range('D31').NumberFormat = "General"
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  -1.72%  '
$ws.Range('E31').NumberFormat = "General"
$ws.Range('E31').Style = "Normal"
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '2.547'
$ws.Range('D32').NumberFormat = "General"
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  +12.15%  '
$ws.Range('E32').NumberFormat = "General"
$ws.Range('E32').Style = "Normal"
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.842.89'
$ws.Range('D33').NumberFormat = "General"
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  -0.68%  '
$ws.Range('E33').NumberFormat = "General"
$ws.Range('E33').Style = "Normal"
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.03015'
$ws.Range('D34').NumberFormat = "General"
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  -5.30%  '
$ws.Range('E34').NumberFormat = "General"
$ws.Range('E34').Style = "Normal"
$ws.Range('B35').NumberFormat = "@"
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('B35').NumberFormat = "General"
$ws.Range('B35').Style = "Normal"
$ws.Range('C35').NumberFormat = "@"
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('C35').NumberFormat = "General"
$ws.Range('C35').Style = "Normal"
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.003'
$ws.Range('D35').NumberFormat = "General"
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  -2.49%  '
$ws.Range('E35').NumberFormat = "General"
$ws.Range('E35').Style = "Normal"
$ws.Range('B36').NumberFormat = "@"
$ws.Range('B36').Value = 'Hedera'
$ws.Range('B36').NumberFormat = "General"
$ws.Range('B36').Style = "Normal"
$ws.Range('C36').NumberFormat = "@"
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('C36').NumberFormat = "General"
$ws.Range('C36').Style = "Normal"
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.08099'
$ws.Range('D36').NumberFormat = "General"
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  -6.20%  '
$ws.Range('E36').NumberFormat = "General"
$ws.Range('E36').Style = "Normal"
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '6.858'
$ws.Range('D37').NumberFormat = "General"
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  -6.30%  '
$ws.Range('E37').NumberFormat = "General"
$ws.Range('E37').Style = "Normal"
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.2762'
$ws.Range('D38').NumberFormat = "General"
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  -3.05%  '
$ws.Range('E38').NumberFormat = "General"
$ws.Range('E38').Style = "Normal"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.09437'
$ws.Range('D39').NumberFormat = "General"
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  +0.06%  '
$ws.Range('E39').NumberFormat = "General"
$ws.Range('E39').Style = "Normal"
$ws.Range('B40').NumberFormat = "@"
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('B40').NumberFormat = "General"
$ws.Range('B40').Style = "Normal"
$ws.Range('C40').NumberFormat = "@"
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('C40').NumberFormat = "General"
$ws.Range('C40').Style = "Normal"
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '10.19'
$ws.Range('D40').NumberFormat = "General"
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  -4.88%  '
$ws.Range('E40').NumberFormat = "General"
$ws.Range('E40').Style = "Normal"
$ws.Range('B41').NumberFormat = "@"
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('B41').NumberFormat = "General"
$ws.Range('B41').Style = "Normal"
$ws.Range('C41').NumberFormat = "@"
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('C41').NumberFormat = "General"
$ws.Range('C41').Style = "Normal"
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.485'
$ws.Range('D41').NumberFormat = "General"
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  +0.51%  '
$ws.Range('E41').NumberFormat = "General"
$ws.Range('E41').Style = "Normal"
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.7779'
$ws.Range('D42').NumberFormat = "General"
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  -6.48%  '
$ws.Range('E42').NumberFormat = "General"
$ws.Range('E42').Style = "Normal"
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '13.36'
$ws.Range('D43').NumberFormat = "General"
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  -5.71%  '
$ws.Range('E43').NumberFormat = "General"
$ws.Range('E43').Style = "Normal"
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '16.33'
$ws.Range('D44').NumberFormat = "General"
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  -7.35%  '
$ws.Range('E44').NumberFormat = "General"
$ws.Range('E44').Style = "Normal"
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  -5.84%  '
$ws.Range('E45').NumberFormat = "General"
$ws.Range('E45').Style = "Normal"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.7019'
$ws.Range('D46').NumberFormat = "General"
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  -5.60%  '
$ws.Range('E46').NumberFormat = "General"
$ws.Range('E46').Style = "Normal"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '4.145'
$ws.Range('D47').NumberFormat = "General"
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  -1.81%  '
$ws.Range('E47').NumberFormat = "General"
$ws.Range('E47').Style = "Normal"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.08587'
$ws.Range('D48').NumberFormat = "General"
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  +2.54%  '
$ws.Range('E48').NumberFormat = "General"
$ws.Range('E48').Style = "Normal"
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.000'
$ws.Range('D49').NumberFormat = "General"
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  +0.00%  '
$ws.Range('E49').NumberFormat = "General"
$ws.Range('E49').Style = "Normal"
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.302'
$ws.Range('D50').NumberFormat = "General"
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  -5.56%  '
$ws.Range('E50').NumberFormat = "General"
$ws.Range('E50').Style = "Normal"
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '136.28'
$ws.Range('D51').NumberFormat = "General"
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  -2.26%  '
$ws.Range('E51').NumberFormat = "General"
$ws.Range('E51').Style = "Normal"
